# Daily TGP price update: shift dates forward one day and update D/E/F/G price columns
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 8
$ws.Range("A8").Value = 45968
$ws.Range("D8").Value = 168.74
$ws.Range("E8").Value = 160.69999999999999
$ws.Range("F8").Value = 170.7
$ws.Range("G8").Value = 160.86000000000001
# Row 9
$ws.Range("A9").Value = 45968
$ws.Range("D9").Value = 168.74
$ws.Range("E9").Value = 160.69999999999999
$ws.Range("F9").Value = 170.7
$ws.Range("G9").Value = 160.86000000000001
# Row 10
$ws.Range("A10").Value = 45968
$ws.Range("D10").Value = 170.96
$ws.Range("E10").Value = 162.94999999999999
$ws.Range("F10").Value = 172.95
$ws.Range("G10").Value = 163.44
# Row 11
$ws.Range("A11").Value = 45967
$ws.Range("D11").Value = 168.28
$ws.Range("E11").Value = 160.72999999999999
$ws.Range("F11").Value = 170.73
$ws.Range("G11").Value = 160.9
# Row 12
$ws.Range("A12").Value = 45967
$ws.Range("D12").Value = 168.28
$ws.Range("E12").Value = 160.72999999999999
$ws.Range("F12").Value = 170.73
$ws.Range("G12").Value = 160.9
# Row 13
$ws.Range("A13").Value = 45967
$ws.Range("D13").Value = 170.11
$ws.Range("E13").Value = 162.91999999999999
$ws.Range("F13").Value = 172.92
$ws.Range("G13").Value = 163.41
# Row 17
$ws.Range("A17").Value = 45968
$ws.Range("D17").Value = 174.44
$ws.Range("E17").Value = 165.92
$ws.Range("F17").Value = 175.92
# Row 18
$ws.Range("A18").Value = 45967
$ws.Range("D18").Value = 173.61
$ws.Range("E18").Value = 165.89
$ws.Range("F18").Value = 175.89
# Row 22
$ws.Range("A22").Value = 45968
$ws.Range("D22").Value = 170.09
$ws.Range("E22").Value = 161.96
$ws.Range("F22").Value = 171.56
$ws.Range("G22").Value = 163.25
# Row 23
$ws.Range("A23").Value = 45968
$ws.Range("D23").Value = 175.73
$ws.Range("E23").Value = 166.67
$ws.Range("F23").Value = 176.67
# Row 24
$ws.Range("A24").Value = 45968
$ws.Range("D24").Value = 175.54
$ws.Range("E24").Value = 166.85
$ws.Range("F24").Value = 176.85
# Row 25
$ws.Range("A25").Value = 45968
$ws.Range("D25").Value = 176.37
$ws.Range("E25").Value = 166.25
$ws.Range("F25").Value = 176.25
$ws.Range("G25").Value = 166.29
# Row 26
$ws.Range("A26").Value = 45968
$ws.Range("D26").Value = 175.1
$ws.Range("E26").Value = 167.81
$ws.Range("F26").Value = 177.81
# Row 27
$ws.Range("A27").Value = 45967
$ws.Range("D27").Value = 169.31
$ws.Range("E27").Value = 161.94999999999999
$ws.Range("F27").Value = 171.55
$ws.Range("G27").Value = 163.22999999999999
# Row 28
$ws.Range("A28").Value = 45967
$ws.Range("D28").Value = 174.88
$ws.Range("E28").Value = 166.65
$ws.Range("F28").Value = 176.65
# Row 29
$ws.Range("A29").Value = 45967
$ws.Range("D29").Value = 174.69
$ws.Range("E29").Value = 166.82
$ws.Range("F29").Value = 176.82
# Row 30
$ws.Range("A30").Value = 45967
$ws.Range("D30").Value = 175.52
$ws.Range("E30").Value = 166.21
$ws.Range("F30").Value = 176.21
$ws.Range("G30").Value = 166.25
# Row 31
$ws.Range("A31").Value = 45967
$ws.Range("D31").Value = 174.26
$ws.Range("E31").Value = 167.77
$ws.Range("F31").Value = 177.77
# Row 35
$ws.Range("A35").Value = 45968
$ws.Range("D35").Value = 169.43
$ws.Range("E35").Value = 160.16
$ws.Range("F35").Value = 169.16
# Row 36
$ws.Range("A36").Value = 45967
$ws.Range("D36").Value = 168.69
$ws.Range("E36").Value = 160.13999999999999
$ws.Range("F36").Value = 169.14
# Row 40
$ws.Range("A40").Value = 45968
$ws.Range("D40").Value = 174.83
$ws.Range("E40").Value = 165.58
$ws.Range("F40").Value = 175.58
# Row 41
$ws.Range("A41").Value = 45968
$ws.Range("D41").Value = 174.53
$ws.Range("E41").Value = 166
$ws.Range("F41").Value = 176
# Row 42
$ws.Range("A42").Value = 45967
$ws.Range("D42").Value = 173.96
$ws.Range("E42").Value = 165.52
$ws.Range("F42").Value = 175.52
# Row 43
$ws.Range("A43").Value = 45967
$ws.Range("D43").Value = 173.67
$ws.Range("E43").Value = 165.94
$ws.Range("F43").Value = 175.95
# Row 47
$ws.Range("A47").Value = 45968
$ws.Range("D47").Value = 168.28
$ws.Range("E47").Value = 161.47999999999999
$ws.Range("F47").Value = 171.48
# Row 48
$ws.Range("A48").Value = 45968
$ws.Range("D48").Value = 168.29
$ws.Range("E48").Value = 161.66999999999999
$ws.Range("F48").Value = 171.67
# Row 49
$ws.Range("A49").Value = 45967
$ws.Range("D49").Value = 168.34
$ws.Range("E49").Value = 161.63999999999999
$ws.Range("F49").Value = 171.64
# Row 50
$ws.Range("A50").Value = 45967
$ws.Range("D50").Value = 168.36
$ws.Range("E50").Value = 161.83000000000001
$ws.Range("F50").Value = 171.83
# Row 54
$ws.Range("A54").Value = 45968
$ws.Range("D54").Value = 185.04
$ws.Range("E54").Value = 175.76
$ws.Range("F54").Value = 185.76
# Row 55
$ws.Range("A55").Value = 45968
$ws.Range("D55").Value = 172.7
$ws.Range("E55").Value = 173.48
$ws.Range("F55").Value = 183.48
# Row 56
$ws.Range("A56").Value = 45968
$ws.Range("D56").Value = 175.19
# Row 57
$ws.Range("A57").Value = 45968
$ws.Range("D57").Value = 174.91
$ws.Range("E57").Value = 167.75
# Row 58
$ws.Range("A58").Value = 45968
$ws.Range("D58").Value = 170.82
$ws.Range("E58").Value = 163.80000000000001
$ws.Range("F58").Value = 173.8
# Row 59
$ws.Range("A59").Value = 45968
$ws.Range("D59").Value = 177.44
$ws.Range("E59").Value = 174.02
# Row 60
$ws.Range("A60").Value = 45967
$ws.Range("D60").Value = 184.18
$ws.Range("E60").Value = 175.67
$ws.Range("F60").Value = 185.67
# Row 61
$ws.Range("A61").Value = 45967
$ws.Range("D61").Value = 171.85
$ws.Range("E61").Value = 173.59
$ws.Range("F61").Value = 183.59
# Row 62
$ws.Range("A62").Value = 45967
$ws.Range("D62").Value = 174.33
# Row 63
$ws.Range("A63").Value = 45967
$ws.Range("D63").Value = 174.09
$ws.Range("E63").Value = 167.85
# Row 64
$ws.Range("A64").Value = 45967
$ws.Range("D64").Value = 170
$ws.Range("E64").Value = 163.91
$ws.Range("F64").Value = 173.91
# Row 65
$ws.Range("A65").Value = 45967
$ws.Range("D65").Value = 176.61
$ws.Range("E65").Value = 173.95
